$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.467.86"

$ws.Range("D3").Value = "1.634.28"
$ws.Range("E3").Value = "  +3.02%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.52"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.83"
$ws.Range("E8").Value = "  +3.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3667"
$ws.Range("E9").Value = "  +2.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.272"
$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08192"
$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.08"
$ws.Range("E13").Value = "  +4.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.641"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001279"
$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.441"
$ws.Range("E16").Value = "  +1.52%  "

$ws.Range("D17").Value = "1.636.93"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.75"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06941"
$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.33"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.568"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "23.471.90"
$ws.Range("E23").Value = "  +1.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.95"
$ws.Range("E24").Value = "  +1.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.102"
$ws.Range("E25").Value = "  +8.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.409"
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.37"
$ws.Range("E27").Value = "  +2.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.04"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.317"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.12"
$ws.Range("E30").Value = "  +3.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.425"
$ws.Range("E31").Value = "  +2.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.799"
$ws.Range("E32").Value = "  +3.14%  "

$ws.Range("D33").Value = "1.815.55"
$ws.Range("E33").Value = "  +2.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9750"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02809"
$ws.Range("E35").Value = "  +5.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.43"
$ws.Range("E36").Value = "  +4.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07425"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.211"
$ws.Range("E38").Value = "  +2.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2537"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08831"
$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.395"
$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7141"
$ws.Range("E42").Value = "  +2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.62"
$ws.Range("E43").Value = "  +4.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.28"
$ws.Range("E44").Value = "  +9.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6606"
$ws.Range("E45").Value = "  +2.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.352"
$ws.Range("E46").Value = "  +3.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.045"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08044"
$ws.Range("E49").Value = "  +2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.94"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("E51").Value = "  +1.37%  "
